$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")
$ws.Range("G1").Value = "TaGs"
